{"js": "// The document has 4 paragraphs:\n//   0: (empty)\n//   1: \"\u0442\u0435\u043c\u0430_\u043e\u0442\u0447\u0435\u0442\u0430\"\n//   2: \"\u0412\u0440\u0435\u043c\u044f \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430: \u0432\u0440\u0435\u043c\u044f\"\n//   3: \"\u0414\u0430\u0442\u0430 \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430: \u0434\u0430\u0442\u0430\"\n//\n// The edit:\n//   1) Right-align paragraphs 2 and 3 (\"\u0412\u0440\u0435\u043c\u044f...\" and \"\u0414\u0430\u0442\u0430...\").\n//   2) Move the auto-managed \"_GoBack\" bookmark so it spans from the very\n//      start of paragraph 2 through the end of the \"\u0434\u0430\u0442\u0430\" run at the end\n//      of paragraph 3 (instead of being an empty bookmark right after\n//      \"\u0434\u0430\u0442\u0430\"), matching where Word leaves it after editing that block.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst timeParagraph = paragraphs.items[2];\nconst dateParagraph = paragraphs.items[3];\n\ntimeParagraph.alignment = Word.Alignment.right;\ndateParagraph.alignment = Word.Alignment.right;\n\n// Remove the existing \"_GoBack\" bookmark (currently collapsed right after\n// the \"\u0434\u0430\u0442\u0430\" text) before re-inserting it over the new span.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nconst goBackStart = timeParagraph.getRange(\"Start\");\nconst goBackEnd = dateParagraph.getRange(\"Content\");\nconst goBackRange = goBackStart.expandTo(goBackEnd);\ngoBackRange.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# The document has 4 paragraphs:\n#   1: (empty)\n#   2: \"\u0442\u0435\u043c\u0430_\u043e\u0442\u0447\u0435\u0442\u0430\"\n#   3: \"\u0412\u0440\u0435\u043c\u044f \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430: \u0432\u0440\u0435\u043c\u044f\"\n#   4: \"\u0414\u0430\u0442\u0430 \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430: \u0434\u0430\u0442\u0430\"\n#\n# The edit:\n#   1) Right-align paragraphs 3 and 4 (\"\u0412\u0440\u0435\u043c\u044f...\" and \"\u0414\u0430\u0442\u0430...\").\n#   2) Move the auto-managed \"_GoBack\" bookmark so it spans from the very\n#      start of paragraph 3 through the end of the \"\u0434\u0430\u0442\u0430\" run at the end\n#      of paragraph 4 (instead of being an empty bookmark right after\n#      \"\u0434\u0430\u0442\u0430\"), matching where Word leaves it after editing that block.\n$d = $word.ActiveDocument\n\n$timeParagraph = $d.Paragraphs(3)\n$dateParagraph = $d.Paragraphs(4)\n\n$timeParagraph.Alignment = \"wdAlignParagraphRight\"\n$dateParagraph.Alignment = \"wdAlignParagraphRight\"\n\n# Remove the existing \"_GoBack\" bookmark (currently collapsed right after\n# the \"\u0434\u0430\u0442\u0430\" text) before re-inserting it over the new span.\n$d.Bookmarks(\"_GoBack\").Delete()\n\n$goBackStart = $timeParagraph.Range.Start\n$goBackEnd = $dateParagraph.Range.End - 1\n$goBackRange = $d.Range($goBackStart, $goBackEnd)\n$d.Bookmarks.Add(\"_GoBack\", $goBackRange)\n"}
